# Inserts one new data row (weekly price report row) above current row 153.
# This pushes the existing rows 153:241 down to 154:242 and the sheet
# dimension grows from A1:R241 to A1:R242 - matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 153 (existing row 153 and everything
# below it shifts down by one row).
$ws.Rows("153:153").Insert()

# Populate the newly inserted row 153 with the new record.
$ws.Range("A153").Value = 7
$ws.Range("B153").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C153").Value = 'Ñuble'
$ws.Range("D153").Value = 44767
$ws.Range("E153").Value = 16
$ws.Range("F153").Value = 100112043
$ws.Range("G153").Value = 'Pepino ensalada'
$ws.Range("H153").Value = 'Sin especificar'
$ws.Range("I153").Value = 'Primera'
$ws.Range("J153").Value = 80
$ws.Range("K153").Value = 19000
$ws.Range("L153").Value = 20000
$ws.Range("M153").Value = 19500
$ws.Range("N153").Value = '$/caja 60 unidades'
$ws.Range("O153").Value = 'Región de Arica y Parinacota'
$ws.Range("P153").Value = 325
$ws.Range("Q153").Value = 60
$ws.Range("R153").Value = 'Hortaliza'
